# Apply "Add data for 2022-04-12" update to the carjacking-by-neighborhood-by-month workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet tab and update the "through" label text (B1 header / shared string).
$ws.Name = "Through 2022-04-04"
$ws.Range("B1").Value = "April 2022 (through April 04)"

# Helper to add an increment to a cell (handles blank cells as 0).
function Add-Count([string]$addr, [double]$amount) {
    $cell = $ws.Range($addr)
    $current = $cell.Value2
    if ($null -eq $current) { $current = 0 }
    $cell.Value = $current + $amount
}

# New carjacking records added for 2022-04-12, incrementing neighborhood/month counts.
Add-Count "B3" 1    # Englewood - April 2022 (through April 04)
Add-Count "R3" 1    # Englewood - April 2018
Add-Count "B5" 1    # Garfield Park - April 2022 (through April 04)
Add-Count "R5" 1    # Garfield Park - April 2018
Add-Count "F6" 1    # Humboldt Park - April 2021
Add-Count "B25" 1   # South Shore - April 2022 (through April 04)
Add-Count "J27" 1   # Washington Heights - April 2020
Add-Count "R38" 1   # Oakland - April 2018
Add-Count "AD39" 1  # Little Village - April 2015
Add-Count "R42" 1   # Avondale - April 2018
Add-Count "B50" 1   # Grand Crossing - April 2022 (through April 04)
Add-Count "J50" 1   # Grand Crossing - April 2020
Add-Count "J51" 1   # Hyde Park - April 2020
Add-Count "F61" 1   # Avalon Park - April 2021
Add-Count "F74" 1   # Little Italy, UIC - April 2021
Add-Count "F91" 1   # West Ridge - April 2021
